# Applies the Jan 5 2023 13:00 symbol-list refresh to the cryptos sheet.
# Re-sets Price (D) / Volume(1h) (E) to the new scraped readings, and
# re-aligns the Coin (B) / Link (C) columns where the source ranking
# reshuffled rows (e.g. rows 9-15 rotated, rows 41/42 swapped).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D/E hold numeric-looking strings ("256.65", "0.57%") that Excel would
# otherwise auto-convert to numbers/percentages on assignment, so those
# cells are forced to Text format first to keep the values as literal
# text, matching the sheet existing text-cell convention.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '256.65'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.57%'
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '27.12'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-1.79%'
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '4.710'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-9.19%'
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05917'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '0.92%'
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.643'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-1.07%'
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8676'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-0.06%'
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9499'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '0.13%'
# Row 9
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1408'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-0.14%'
# Row 10
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.03830'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '11.29%'
# Row 11
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07096'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-0.99%'
# Row 12
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03217'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '1.19%'
# Row 13
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09263'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.33%'
# Row 14
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001534'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.89%'
# Row 15
$ws.Range("B15").Value = 'One'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0006029'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.81%'
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006022'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.35%'
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.514'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.43%'
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.194'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-1.03%'
# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.66%'
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.3134'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-1.24%'
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1283'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-0.75%'
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.880'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '9.36%'
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04233'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '1.35%'
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001223'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-0.30%'
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004297'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-10.30%'
# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-0.01%'
# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '2.42%'
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03820'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '0.34%'
# Row 41
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1102'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-0.13%'
# Row 42
$ws.Range("B42").Value = 'KickToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.003954'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-30.02%'
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002420'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '0.00%'
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01151'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '17.38%'
# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '2.49%'
# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '0.02%'
# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-19.49%'
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002281'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '7.22%'
# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '0.02%'
# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.02%'
